$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.924.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.536.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.28%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.66%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.579'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.535.99'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.61'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.35%  '

$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.991.45'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.876.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000141'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.541.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '334.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.94%  '

$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.59'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.15%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("B28").Value = 'SuiNetwork'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.81%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0812'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '177.59'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '412.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.396'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.43%  '

$ws.Range("B38").Value = 'USDe'
$ws.Range("C38").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.08%  '

$ws.Range("E41").Value = '  +0.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.88%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.47'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.604'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0963'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0516'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0234'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.36'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.45%  '
